$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.855
$ws.Range("D18").Value = -8.292000000000002
$ws.Range("B21").Value = 9.019
$ws.Range("B23").Value = 8.843999999999999
$ws.Range("C24").Value = -12.015
$ws.Range("B25").Value = 5.922000000000001
$ws.Range("C28").Value = -12.938
$ws.Range("C36").Value = -12.732
$ws.Range("C45").Value = -13.361
$ws.Range("C48").Value = -11.175
$ws.Range("C49").Value = -12.754
$ws.Range("D51").Value = -8.272000000000002
$ws.Range("C52").Value = -10.931
$ws.Range("B53").Value = 5.389
$ws.Range("C53").Value = -10.494
$ws.Range("C54").Value = -12.907
$ws.Range("D55").Value = -8.177000000000001
$ws.Range("B57").Value = 5.061999999999999
$ws.Range("B59").Value = 4.708
$ws.Range("D64").Value = -7.454000000000001
$ws.Range("B69").Value = 5.726
$ws.Range("C70").Value = -11.3
$ws.Range("B79").Value = 5.577
$ws.Range("D80").Value = -7.973999999999999
$ws.Range("B83").Value = 5.165999999999999
$ws.Range("C86").Value = -13.477
$ws.Range("C87").Value = -13.621
$ws.Range("D92").Value = -6.842000000000001
$ws.Range("B93").Value = 5.131000000000001
$ws.Range("D94").Value = -7.101000000000001
$ws.Range("D96").Value = -7.414
$ws.Range("C101").Value = -12.26
